$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 379, pushing existing rows 379:445 down to 380:446.
$ws.Rows("379:379").Insert()

# Populate the new row 379 with a fresh weekly data point (same dimensions
# as the old row 379, but with updated Fecha and Volumen).
$ws.Cells.Item(379, 1).Value = 10
$ws.Cells.Item(379, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(379, 3).Value = "La Araucanía"
$ws.Cells.Item(379, 4).Value = 44951
$ws.Cells.Item(379, 5).Value = 9
$ws.Cells.Item(379, 6).Value = 100112009
$ws.Cells.Item(379, 7).Value = "Acelga"
$ws.Cells.Item(379, 8).Value = "Sin especificar"
$ws.Cells.Item(379, 9).Value = "Primera"
$ws.Cells.Item(379, 10).Value = 35
$ws.Cells.Item(379, 11).Value = 8000
$ws.Cells.Item(379, 12).Value = 8000
$ws.Cells.Item(379, 13).Value = 8000
$ws.Cells.Item(379, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(379, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(379, 16).Value = 667
$ws.Cells.Item(379, 17).Value = 12
$ws.Cells.Item(379, 18).Value = "Hortaliza"
